# Apply updated crypto price/volume figures to sheet1 (cryptos list refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells below are plain-text table cells (prices use "." as thousands
# separators, volumes carry a leading/trailing double space). Force the
# NumberFormat to Text ("@") before assignment so Excel does not silently
# reinterpret numeric-looking strings (e.g. "107.98") as actual numbers,
# then clear the format back off so no stray style index is left behind.
$cellUpdates = @{
    'D2' = '45.398.24'
    'E2' = '  +0.01%  '
    'D3' = '2.369.73'
    'E3' = '  -0.44%  '
    'E4' = '  -0.01%  '
    'E5' = '  -1.57%  '
    'D6' = '107.98'
    'E6' = '  -3.56%  '
    'E7' = '  -1.30%  '
    'E8' = '  -0.01%  '
    'E9' = '  -3.52%  '
    'D10' = '40.75'
    'E10' = '  -3.74%  '
    'D11' = '0.0918'
    'E11' = '  -1.32%  '
    'E12' = '  -2.74%  '
    'E13' = '  +0.72%  '
    'E14' = '  -3.44%  '
    'D15' = '2.731.33'
    'E15' = '  -0.46%  '
    'D16' = '15.32'
    'E16' = '  -3.25%  '
    'D17' = '2.372.22'
    'E17' = '  -0.74%  '
    'D18' = '45.469.89'
    'E18' = '  +0.23%  '
    'D19' = '13.86'
    'E19' = '  +5.75%  '
    'E20' = '  -1.63%  '
    'D21' = '7.24'
    'E21' = '  -5.10%  '
    'D22' = '73.31'
    'E22' = '  -2.58%  '
    'E23' = '  -0.77%  '
    'D24' = '259.91'
    'E24' = '  -3.45%  '
    'D25' = '2.39'
    'E25' = '  +2.04%  '
    'E26' = '  -0.06%  '
    'D27' = '11.16'
    'D28' = '7.19'
    'E28' = '  -6.72%  '
    'E29' = '  -1.45%  '
    'D30' = '0.0978'
    'E30' = '  +4.62%  '
    'D31' = '22.34'
    'E31' = '  -2.46%  '
    'D32' = '37.17'
    'E32' = '  -5.44%  '
    'D33' = '166.70'
    'E33' = '  -1.75%  '
    'D34' = '2.96'
    'E34' = '  -0.29%  '
    'E35' = '  -2.18%  '
    'E36' = '  +0.71%  '
    'D37' = '4.72'
    'E37' = '  -2.53%  '
    'D38' = '4.00'
    'E38' = '  +1.46%  '
    'D39' = '1.90'
    'E39' = '  +8.57%  '
    'E40' = '  -3.64%  '
    'E41' = '  -3.74%  '
    'D42' = '99.59'
    'E42' = '  -5.93%  '
    'D43' = '69.80'
    'E43' = '  -3.28%  '
    'E44' = '  -5.47%  '
    'E45' = '  -0.09%  '
    'D46' = '12.66'
    'E46' = '  -8.73%  '
    'D47' = '1.820.57'
    'E47' = '  +9.81%  '
    'D48' = '84.78'
    'E48' = '  +5.91%  '
    'E49' = '  +1.27%  '
    'D50' = '9.26'
    'E50' = '  +1.70%  '
    'D51' = '110.80'
    'E51' = '  -6.99%  '
}

foreach ($addr in $cellUpdates.Keys) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $cellUpdates[$addr]
    $rng.ClearFormats()
}
